$d = $word.ActiveDocument

# --- 1. Fix typo "unrelated to the robot managed by the robot" -> "...the battery managed by the robot"
$d.Content.Find.Execute("unrelated to the robot managed by the robot", $true, $false, $false, $false, $false,
                         $true, 1, $false, "unrelated to the battery managed by the robot", 2) | Out-Null

Write-Output "done"
